{"js": "// Replace the division-problem text in the worksheet table.\n// The table is 20 rows x 5 cols; only every 4th row (0, 4, 8, 12, 16) has\n// visible content, one \"NN\u00f7N=\" expression per cell. We update each of\n// those 25 cells in row-major (reading) order to the new values, using\n// row/column addressing instead of a global text search-and-replace so\n// that duplicate values (e.g. \"93\u00f78=\" or \"77\u00f78=\" occurring more than\n// once in the original) each get their own, independent replacement,\n// exactly matching the order the values appear in the document.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row index -> ordered list of [oldText, newText] for that row's 5 cells.\nconst rowUpdates = {\n  0: [\n    [\"77\u00f72=\", \"63\u00f78=\"],\n    [\"57\u00f76=\", \"87\u00f79=\"],\n    [\"77\u00f77=\", \"63\u00f78=\"],\n    [\"82\u00f72=\", \"89\u00f79=\"],\n    [\"28\u00f74=\", \"29\u00f72=\"],\n  ],\n  4: [\n    [\"11\u00f76=\", \"23\u00f73=\"],\n    [\"13\u00f75=\", \"12\u00f75=\"],\n    [\"34\u00f76=\", \"12\u00f77=\"],\n    [\"82\u00f78=\", \"98\u00f72=\"],\n    [\"83\u00f74=\", \"99\u00f72=\"],\n  ],\n  8: [\n    [\"93\u00f78=\", \"92\u00f79=\"],\n    [\"36\u00f76=\", \"17\u00f75=\"],\n    [\"37\u00f79=\", \"37\u00f75=\"],\n    [\"91\u00f76=\", \"64\u00f73=\"],\n    [\"93\u00f78=\", \"25\u00f77=\"],\n  ],\n  12: [\n    [\"39\u00f74=\", \"91\u00f73=\"],\n    [\"57\u00f78=\", \"30\u00f76=\"],\n    [\"77\u00f78=\", \"67\u00f72=\"],\n    [\"14\u00f78=\", \"35\u00f78=\"],\n    [\"41\u00f75=\", \"42\u00f79=\"],\n  ],\n  16: [\n    [\"36\u00f75=\", \"68\u00f79=\"],\n    [\"99\u00f73=\", \"88\u00f73=\"],\n    [\"86\u00f74=\", \"75\u00f76=\"],\n    [\"77\u00f78=\", \"40\u00f75=\"],\n    [\"75\u00f75=\", \"94\u00f79=\"],\n  ],\n};\n\nfor (const rowIndex of Object.keys(rowUpdates)) {\n  const row = parseInt(rowIndex, 10);\n  const cells = rowUpdates[rowIndex];\n  for (let col = 0; col < cells.length; col++) {\n    const [, newText] = cells[col];\n    const cell = table.getCell(row, col);\n    // Setting `.value` rewrites the cell's text while leaving the\n    // existing run/paragraph formatting (font, size, alignment) intact.\n    cell.value = newText;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the division-problem text in the worksheet table.\n# The table is 20 rows x 5 cols; only every 4th row (Word's 1-based rows\n# 1, 5, 9, 13, 17) has visible content, one \"NN\u00f7N=\" expression per cell.\n# We address each of those 25 cells directly by (row, column) rather than\n# doing a global text search-and-replace, so that duplicate values\n# (e.g. \"93\u00f78=\" or \"77\u00f78=\" occurring more than once in the original) are\n# each replaced independently with the correct new value, matching the\n# order the values appear in the document.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Each element: row, col, expected old text, new text.\n$updates = @(\n    @(1, 1, \"77\u00f72=\", \"63\u00f78=\"),\n    @(1, 2, \"57\u00f76=\", \"87\u00f79=\"),\n    @(1, 3, \"77\u00f77=\", \"63\u00f78=\"),\n    @(1, 4, \"82\u00f72=\", \"89\u00f79=\"),\n    @(1, 5, \"28\u00f74=\", \"29\u00f72=\"),\n\n    @(5, 1, \"11\u00f76=\", \"23\u00f73=\"),\n    @(5, 2, \"13\u00f75=\", \"12\u00f75=\"),\n    @(5, 3, \"34\u00f76=\", \"12\u00f77=\"),\n    @(5, 4, \"82\u00f78=\", \"98\u00f72=\"),\n    @(5, 5, \"83\u00f74=\", \"99\u00f72=\"),\n\n    @(9, 1, \"93\u00f78=\", \"92\u00f79=\"),\n    @(9, 2, \"36\u00f76=\", \"17\u00f75=\"),\n    @(9, 3, \"37\u00f79=\", \"37\u00f75=\"),\n    @(9, 4, \"91\u00f76=\", \"64\u00f73=\"),\n    @(9, 5, \"93\u00f78=\", \"25\u00f77=\"),\n\n    @(13, 1, \"39\u00f74=\", \"91\u00f73=\"),\n    @(13, 2, \"57\u00f78=\", \"30\u00f76=\"),\n    @(13, 3, \"77\u00f78=\", \"67\u00f72=\"),\n    @(13, 4, \"14\u00f78=\", \"35\u00f78=\"),\n    @(13, 5, \"41\u00f75=\", \"42\u00f79=\"),\n\n    @(17, 1, \"36\u00f75=\", \"68\u00f79=\"),\n    @(17, 2, \"99\u00f73=\", \"88\u00f73=\"),\n    @(17, 3, \"86\u00f74=\", \"75\u00f76=\"),\n    @(17, 4, \"77\u00f78=\", \"40\u00f75=\"),\n    @(17, 5, \"75\u00f75=\", \"94\u00f79=\")\n)\n\nforeach ($u in $updates) {\n    $row = $u[0]\n    $col = $u[1]\n    $newText = $u[3]\n    $cell = $t.Cell($row, $col)\n    # Assigning to the cell's Range.Text replaces the visible characters\n    # (Word keeps the trailing cell-mark) while leaving the existing\n    # run/paragraph formatting (font, size, alignment) untouched.\n    $cell.Range.Text = $newText\n}\n"}
